$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp title
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 12:22"

# Row 36: 'Rumania' -> 'Rumania'
$ws.Range("B36").Value = 12240
$ws.Range("C36").Value = 262
$ws.Range("D36").Value = 4017
$ws.Range("E36").Value = 7528
$ws.Range("F36").Value = 221

# Row 54: 'Finlandia' -> 'Finlandia'
$ws.Range("B54").Value = 4995
$ws.Range("C54").Value = 89
$ws.Range("E54").Value = 1989

# Row 55: 'Marruecos' -> 'Marruecos'
$ws.Range("B55").Value = 4359
$ws.Range("C55").Value = 38
$ws.Range("D55").Value = 969
$ws.Range("E55").Value = 3222

# Row 76: 'Nigeria' -> 'Bosnia y Herzegovina'
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1757
$ws.Range("C76").Value = 80
$ws.Range("D76").Value = 727
$ws.Range("E76").Value = 961
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 69

# Row 77: 'Estonia' -> 'Nigeria'
$ws.Range("A77").Value = "Nigeria"
$ws.Range("B77").Value = 1728
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 307
$ws.Range("E77").Value = 1370
$ws.Range("F77").Value = 2
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 51

# Row 78: 'Bosnia y Herzegovina' -> 'Estonia'
$ws.Range("A78").Value = "Estonia"
$ws.Range("B78").Value = 1689
$ws.Range("C78").Value = 23
$ws.Range("D78").Value = 249
$ws.Range("E78").Value = 1388
$ws.Range("F78").Value = 9
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 52

# Row 90: 'Hong Kong' -> 'Hong Kong'
$ws.Range("D90").Value = 846
$ws.Range("E90").Value = 188

# Row 170: 'Macao' -> 'Macao'
$ws.Range("D170").Value = 35
$ws.Range("E170").Value = 10

# Row 196: 'Nicaragua' -> 'Santo Tome y Principe'
$ws.Range("A196").Value = "Santo Tome y Principe"
$ws.Range("B196").Value = 14
$ws.Range("C196").Value = 6
$ws.Range("D196").Value = 4
$ws.Range("E196").Value = 10
$ws.Range("H196").Value = 0

# Row 197: 'Islas Malvinas' -> 'Nicaragua'
$ws.Range("A197").Value = "Nicaragua"
$ws.Range("D197").Value = 7
$ws.Range("E197").Value = 3
$ws.Range("H197").Value = 3

# Row 198: 'Islas Turcas y Caicos' -> 'Islas Malvinas'
$ws.Range("A198").Value = "Islas Malvinas"
$ws.Range("B198").Value = 13
$ws.Range("D198").Value = 11
$ws.Range("E198").Value = 2
$ws.Range("H198").Value = 0

# Row 199: 'Montserrat' -> 'Islas Turcas y Caicos'
$ws.Range("A199").Value = "Islas Turcas y Caicos"
$ws.Range("B199").Value = 12
$ws.Range("D199").Value = 5
$ws.Range("E199").Value = 6
$ws.Range("F199").Value = 0

# Row 200: 'Burundi' -> 'Montserrat'
$ws.Range("A200").Value = "Montserrat"
$ws.Range("D200").Value = 2
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = 1

# Row 201: 'Seychelles' -> 'Burundi'
$ws.Range("A201").Value = "Burundi"
$ws.Range("D201").Value = 4
$ws.Range("E201").Value = 6
$ws.Range("H201").Value = 1

# Row 202: 'Gambia' -> 'Seychelles'
$ws.Range("A202").Value = "Seychelles"
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 6
$ws.Range("E202").Value = 5
$ws.Range("H202").Value = 0

# Row 203: 'Groenlandia' -> 'Gambia'
$ws.Range("A203").Value = "Gambia"
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 8
$ws.Range("E203").Value = 2
$ws.Range("H203").Value = 1

# Row 204: 'Santa Sede' -> 'Groenlandia'
$ws.Range("A204").Value = "Groenlandia"
$ws.Range("B204").Value = 11
$ws.Range("D204").Value = 11
$ws.Range("E204").Value = 0

# Row 205: 'Santo Tome y Principe' -> 'Santa Sede'
$ws.Range("A205").Value = "Santa Sede"
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 8
